$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 642.36365
$ws.Range("I58").Value = 285.66666
$ws.Range("K58").Value = 856.9999799999999
$ws.Range("M58").Value = -706.9999799999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3254.5715
$ws.Range("I98").Value = 3187.6365
$ws.Range("K98").Value = 3187.6365
$ws.Range("M98").Value = -1689.6365

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3254.5715
$ws.Range("I122").Value = 3187.6365
$ws.Range("K122").Value = 9562.9095
$ws.Range("M122").Value = -7112.9095

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1903.7587
$ws.Range("I137").Value = 1836.1177
$ws.Range("J137").Value = 1999.5834
$ws.Range("K137").Value = 5508.3531
$ws.Range("L137").Value = 5998.7502
$ws.Range("M137").Value = -2958.3531
$ws.Range("N137").Value = -11098.7502

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 922.8
$ws.Range("I141").Value = 922.8
$ws.Range("K141").Value = 2768.4
$ws.Range("M141").Value = 2411.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7218.4443
$ws.Range("I32").Value = 3811.4333
$ws.Range("J32").Value = 24253.5
$ws.Range("K32").Value = 3811.4333
$ws.Range("L32").Value = 24253.5
$ws.Range("M32").Value = -3524.4333
$ws.Range("N32").Value = -24827.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2080
$ws.Range("I122").Value = 1607.75
$ws.Range("K122").Value = 4823.25
$ws.Range("M122").Value = -2373.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 28034.5
$ws.Range("J125").Value = 28034.5
$ws.Range("L125").Value = 28034.5
$ws.Range("N125").Value = -37874.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2132.6428
$ws.Range("I132").Value = 1867.5143
$ws.Range("K132").Value = 5602.5429
$ws.Range("M132").Value = -3072.5429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3271.2964
$ws.Range("I20").Value = 2759.0605
$ws.Range("J20").Value = 4076.238
$ws.Range("K20").Value = 2759.0605
$ws.Range("L20").Value = 4076.238
$ws.Range("M20").Value = -2512.0605
$ws.Range("N20").Value = -4570.237999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3050.963
$ws.Range("I58").Value = 2873.647
$ws.Range("J58").Value = 3352.4
$ws.Range("K58").Value = 2873.647
$ws.Range("L58").Value = 3352.4
$ws.Range("M58").Value = -2670.647
$ws.Range("N58").Value = -3758.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 17986
$ws.Range("J60").Value = 18180.6
$ws.Range("L60").Value = 18180.6
$ws.Range("N60").Value = -19202.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 14099.878
$ws.Range("I134").Value = 4724.8887
$ws.Range("J134").Value = 81599.8
$ws.Range("K134").Value = 14174.6661
$ws.Range("L134").Value = 244799.4
$ws.Range("M134").Value = -11639.6661
$ws.Range("N134").Value = -249869.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3050.963
$ws.Range("I136").Value = 2873.647
$ws.Range("J136").Value = 3352.4
$ws.Range("K136").Value = 8620.940999999999
$ws.Range("L136").Value = 10057.2
$ws.Range("M136").Value = -6070.940999999999
$ws.Range("N136").Value = -15157.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 435.5
$ws.Range("I11").Value = 435.5
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1306.5
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1166.5
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 6260
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 6575
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 19725
$ws.Range("M70").Value = -14685
$ws.Range("N70").Value = -20355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 6260
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 6575
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 19725
$ws.Range("M73").Value = -13908
$ws.Range("N73").Value = -21909

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 660
$ws.Range("I94").Value = 660
$ws.Range("K94").Value = 1980
$ws.Range("M94").Value = -1304

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 7166.5
$ws.Range("I104").Value = 4333.3335
$ws.Range("K104").Value = 13000.0005
$ws.Range("M104").Value = -10379.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 11800
$ws.Range("J106").Value = 12500
$ws.Range("L106").Value = 37500
$ws.Range("N106").Value = -39392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 24831.818
$ws.Range("J131").Value = 2008.75
$ws.Range("L131").Value = 6026.25
$ws.Range("N131").Value = -16106.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 101121
$ws.Range("J140").Value = 101121
$ws.Range("L140").Value = 101121
$ws.Range("N140").Value = -111481

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 119246.5
$ws.Range("J141").Value = 119246.5
$ws.Range("L141").Value = 119246.5
$ws.Range("N141").Value = -129606.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 7571.2856
$ws.Range("J20").Value = 5999.5
$ws.Range("L20").Value = 5999.5
$ws.Range("N20").Value = -6451.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3455.25
$ws.Range("I68").Value = 3599.7
$ws.Range("J68").Value = 2733
$ws.Range("K68").Value = 3599.7
$ws.Range("L68").Value = 2733
$ws.Range("M68").Value = -2850.7
$ws.Range("N68").Value = -4231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3455.25
$ws.Range("I71").Value = 3599.7
$ws.Range("J71").Value = 2733
$ws.Range("K71").Value = 17998.5
$ws.Range("L71").Value = 13665
$ws.Range("M71").Value = -14254.5
$ws.Range("N71").Value = -21153

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10939.733
$ws.Range("I81").Value = 51749.5
$ws.Range("J81").Value = 4661.3076
$ws.Range("K81").Value = 103499
$ws.Range("L81").Value = 9322.6152
$ws.Range("M81").Value = -102438
$ws.Range("N81").Value = -11444.6152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 10939.733
$ws.Range("I84").Value = 51749.5
$ws.Range("J84").Value = 4661.3076
$ws.Range("K84").Value = 517495
$ws.Range("L84").Value = 46613.076
$ws.Range("M84").Value = -512191
$ws.Range("N84").Value = -57221.076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 26000
$ws.Range("J104").Value = 26000
$ws.Range("L104").Value = 26000
$ws.Range("N104").Value = -32988

Write-Host "Edit complete"
